$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 ("I0") and J1 ("IF"), copying the existing
# header formatting (bold font, borders, centered alignment) from H1
# so they share the same cell style as the other header cells.
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("I1").Value = "I0"

$ws.Range("H1").Copy()
$ws.Range("J1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("J1").Value = "IF"

# Populate the new I and J columns with their data values (rows 2-29)
$iValues = @(8,8,5,9,1,1,1,1,7,8,7,5,8,3,1,1,1,8,1,1,7,1,1,1,1,1,1,1)
$jValues = @(9,9,7,9,6,6,5,5,7,8,8,6,8,7,7,5,3,8,5,6,7,4,5,5,6,5,3,2)

for ($i = 0; $i -lt 28; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$i]
    $ws.Cells.Item($row, 10).Value = $jValues[$i]
}
